$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("PDL Date") holds date-like text such as "03052025".
# Update every data row (2-303) to the new placeholder date "99999999".
# Force a text number format first so Excel keeps storing these as text
# (matching the original inline-string text cells) instead of coercing
# the all-digit value into a numeric cell.
$dateRange = $ws.Range("B2:B303")
$dateRange.NumberFormat = "@"
$dateRange.Value = "99999999"

# Fix the D179/D180 "Orencia" drug description: the embedded line break
# between the comma and "CLICKJECT" becomes a single space.
$ws.Cells.Item(179, 4).Value = "ORENCIA 50mg/0.4ml, 87.5mg/0.7ml, 125MG/ML, CLICKJECT"
$ws.Cells.Item(180, 4).Value = "ORENCIA 50mg/0.4ml, 87.5mg/0.7ml, 125MG/ML, CLICKJECT"
